$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (Home) target depth data update for Week 17
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 379
$wsOff.Range("C2").Value = 294
$wsOff.Range("D2").Value = 98
$wsOff.Range("E2").Value = 52

# Sheet "DEF" - row 2 (Home) target depth data update for Week 17
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 515
$wsDef.Range("C2").Value = 362
$wsDef.Range("D2").Value = 122
$wsDef.Range("E2").Value = 55
